$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update difficulty ("Difficulté (F,PF)") values in column G for rows 8 and 9
$ws.Range("G8").Value = "PF"
$ws.Range("G9").Value = "PF"

# Update the active selection to match the saved view state
$ws.Activate()
$ws.Range("F4").Select()
